$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.343.48"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.933.82"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'357.73"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'110.46"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.632"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'39.17"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'19.59"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "'7.83"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "3.398.03"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "2.945.78"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "'0.989"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "52.315.96"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "'14.03"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "0.0₃0986"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'70.61"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'269.38"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  +6.65%  "
$ws.Range("E27").Value = "  +15.87%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'0.107"
$ws.Range("E30").Value = "  +7.71%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'37.71"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'6.26"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "'0.0444"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'3.20"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("D39").Value = "'18.33"
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'2.00"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'23.08"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "'119.72"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.48"
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.47"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").Value = "2.133.72"
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "'0.931"
$ws.Range("E51").Value = "  -3.98%  "
